# Insert a new data row at row 139 (pushing existing rows 139-262 down to 140-263)
# and populate it with the new weekly price record for Pepino dulce.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(139).Insert()

$ws.Cells.Item(139, 1).Value = 10
$ws.Cells.Item(139, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(139, 3).Value = "La Araucanía"
$ws.Cells.Item(139, 4).Value = 44778
$ws.Cells.Item(139, 5).Value = 9
$ws.Cells.Item(139, 6).Value = 100112043
$ws.Cells.Item(139, 7).Value = "Pepino dulce"
$ws.Cells.Item(139, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(139, 9).Value = "Primera"
$ws.Cells.Item(139, 10).Value = 95
$ws.Cells.Item(139, 11).Value = 20000
$ws.Cells.Item(139, 12).Value = 20000
$ws.Cells.Item(139, 13).Value = 20000
$ws.Cells.Item(139, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(139, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(139, 16).Value = 1111
$ws.Cells.Item(139, 17).Value = 18
$ws.Cells.Item(139, 18).Value = "Hortaliza"

# Preserve the date style (s="2") that the row above/below use, in case Insert
# did not already copy it through.
$ws.Cells.Item(139, 4).NumberFormat = $ws.Cells.Item(140, 4).NumberFormat
